$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for rows 2-11 (Q0..Q9), columns B(ME) C(MAE) D(MSE) E(RMSE) F(SE) G(N)

$ws.Range("B2").Value = 0.08102716014508599
$ws.Range("C2").Value = 0.1849562581409688
$ws.Range("D2").Value = 0.05080724786734196
$ws.Range("E2").Value = 0.225404631423895
$ws.Range("F2").Value = 0.2177199557413915
$ws.Range("G2").Value = 15

$ws.Range("B3").Value = 0.3284388860093228
$ws.Range("C3").Value = 0.3436510707338022
$ws.Range("D3").Value = 0.1900396029840986
$ws.Range("E3").Value = 0.4359353197254136
$ws.Range("F3").Value = 0.2974694575110258
$ws.Range("G3").Value = 14

$ws.Range("B4").Value = 0.5008085514463395
$ws.Range("C4").Value = 0.5008085514463395
$ws.Range("D4").Value = 0.4390583335671841
$ws.Range("E4").Value = 0.6626147701094386
$ws.Range("F4").Value = 0.4515933521755904
$ws.Range("G4").Value = 13

$ws.Range("B5").Value = 0.6271289531314115
$ws.Range("C5").Value = 0.6271289531314115
$ws.Range("D5").Value = 0.6001355647093024
$ws.Range("E5").Value = 0.77468417094278
$ws.Range("F5").Value = 0.4750251754326699
$ws.Range("G5").Value = 12

$ws.Range("B6").Value = 0.6196096897755812
$ws.Range("C6").Value = 0.6196096897755812
$ws.Range("D6").Value = 0.4847206317237406
$ws.Range("E6").Value = 0.6962188102340676
$ws.Range("F6").Value = 0.3329938595018584
$ws.Range("G6").Value = 11

$ws.Range("B7").Value = 0.4878350037768636
$ws.Range("C7").Value = 0.4878350037768636
$ws.Range("D7").Value = 0.2918424992789185
$ws.Range("E7").Value = 0.5402244897067501
$ws.Range("F7").Value = 0.2446301252661204
$ws.Range("G7").Value = 10

$ws.Range("B8").Value = 0.408887426949955
$ws.Range("C8").Value = 0.408887426949955
$ws.Range("D8").Value = 0.2151763681766355
$ws.Range("E8").Value = 0.4638710684841592
$ws.Range("F8").Value = 0.2323485964907919
$ws.Range("G8").Value = 9

$ws.Range("B9").Value = 0.399669795446622
$ws.Range("C9").Value = 0.399669795446622
$ws.Range("D9").Value = 0.2150887163119547
$ws.Range("E9").Value = 0.4637765801676004
$ws.Range("F9").Value = 0.2577272300389153
$ws.Range("G9").Value = 6

$ws.Range("B10").Value = 0.2327354214815493
$ws.Range("C10").Value = 0.2327354214815493
$ws.Range("D10").Value = 0.07574075157082642
$ws.Range("E10").Value = 0.2752103769315874
$ws.Range("F10").Value = 0.1798956996093793
$ws.Range("G10").Value = 3

$ws.Range("B11").Value = 0.2378331735378737
$ws.Range("C11").Value = 0.2378331735378737
$ws.Range("D11").Value = 0.05656461843509633
$ws.Range("E11").Value = 0.2378331735378737
$ws.Range("F11").ClearContents()
$ws.Range("G11").Value = 1
